# Generate Report for Handoff
#
# The localization-status report is regenerated: the "In Translation" status
# becomes "Ready for handoff" everywhere it appears (Overview zh-cn/de-de
# status columns plus the per-language Status column), the "Latest HO Xliff
# Generate Date" / "Latest Handoff Datetime" timestamps are refreshed to the
# new run time, and the Status column is widened slightly to fit the new
# (longer) text on the Overview, zh-cn and de-de sheets.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

# --- Status text: "In Translation" -> "Ready for handoff" ---------------
$overview.Range("E2").Value = "Ready for handoff"
$overview.Range("F2").Value = "Ready for handoff"
$zhcn.Range("C2").Value     = "Ready for handoff"
$dede.Range("C2").Value     = "Ready for handoff"

# --- Refreshed timestamps -------------------------------------------------
$overview.Range("G2").Value = "2016-08-18 04:58:43"
$dede.Range("H2").Value     = "2016-08-18 04:58:43"
$zhcn.Range("H2").Value     = "2016-08-18 04:58:37"

# --- Widen the Status columns so the new text fits -------------------------
# (Target width ~17.22 chars; Excel snaps ColumnWidth to whole-pixel steps at
# the active font's Normal-style max-digit-width, so 16.33 is the input that
# lands on the nearest attainable grid value.)
$overview.Columns.Item(5).ColumnWidth = 16.33
$overview.Columns.Item(6).ColumnWidth = 16.33
$zhcn.Columns.Item(3).ColumnWidth     = 16.33
$dede.Columns.Item(3).ColumnWidth     = 16.33
